$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A4 was previously stored as an inline string "76442781"; it should be a
# real number (matches A2/A3 which already hold numeric phone values).
$ws.Range("A4").Value = 76442781

# Add new row 5 for payment 79174463 (Cash) 2025-08-20T09:36:16
# Phone numbers are kept as text even though they look numeric (consistent
# with how A3 -> A4 originally stored this same phone number as text), so
# force text formatting before assigning, then clear the number-format
# override so no extra style is left behind on the cell.
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "79174463"
$ws.Range("A5").ClearFormats()

$ws.Range("B5").Value = 40
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 40
$ws.Range("G5").Value = "Cash"
$ws.Range("H5").Value = "2025-08-20T09:36:16"
